$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Col5a1", "Sdc3", "ECs", 3, 1, 1.067323333333333, 3.20197, 0.002262082342179458, 0.002262082342179457, 3, 1, 24.44779533333333, 73.343386, 0.1890645123346783, 0.1890645123346783, 26.09370240782444, 234.84332167042, 0.0004276794948850461, 0.0004276794948850461),
    @("ECs", "Col5a1", "Sdc3", "FAPs", 3, 1, 1.067323333333333, 3.20197, 0.002262082342179458, 0.002262082342179457, 3, 1, 18.11074233333333, 54.332227, 0.1400575643155068, 0.1400575643155068, 19.33001787635444, 173.97016088719, 0.0003168217431267718, 0.0003168217431267717),
    @("ECs", "Col5a1", "Sdc3", "M2", 3, 1, 1.067323333333333, 3.20197, 0.002262082342179458, 0.002262082342179457, 3, 1, 80.684877, 242.054631, 0.6239682030546764, 0.6239682030546765, 86.11685186922999, 775.0516668230699, 0.00141146745421143, 0.00141146745421143),
    @("ECs", "Col5a1", "Sdc3", "sCs", 3, 1, 1.067323333333333, 3.20197, 0.002262082342179458, 0.002262082342179457, 3, 1, 6.065862, 18.197586, 0.0469097202951384, 0.04690972029513841, 6.474236049379999, 58.26812444442, 0.0001061136499562099, 0.0001061136499562099),
    @("FAPs", "Col5a1", "Sdc3", "ECs", 3, 1, 416.538737, 1249.616211, 0.8828111335847306, 0.8828111335847305, 3, 1, 24.44779533333333, 73.343386, 0.1890645123346783, 0.1890645123346783, 10183.45379058116, 91651.08411523043, 0.1669082564548216, 0.1669082564548216),
    @("FAPs", "Col5a1", "Sdc3", "FAPs", 3, 1, 416.538737, 1249.616211, 0.8828111335847306, 0.8828111335847305, 3, 1, 18.11074233333333, 54.332227, 0.1400575643155068, 0.1400575643155068, 7543.825737659101, 67894.4316389319, 0.1236443771204889, 0.1236443771204889),
    @("FAPs", "Col5a1", "Sdc3", "M2", 3, 1, 416.538737, 1249.616211, 0.8828111335847306, 0.8828111335847305, 3, 1, 80.684877, 242.054631, 0.6239682030546764, 0.6239682030546765, 33608.37676058035, 302475.3908452232, 0.5508460766595262, 0.5508460766595262),
    @("FAPs", "Col5a1", "Sdc3", "sCs", 3, 1, 416.538737, 1249.616211, 0.8828111335847306, 0.8828111335847305, 3, 1, 6.065862, 18.197586, 0.0469097202951384, 0.04690972029513841, 2526.666496296294, 22739.99846666665, 0.04141242334989378, 0.04141242334989378),
    @("M2", "Col5a1", "Sdc3", "ECs", 2, 0.6666666666666666, 0.1415573333333333, 0.424672, 0.0003000162501266517, 0.0003000162501266516, 3, 1, 24.44779533333333, 73.343386, 0.1890645123346783, 0.1890645123346783, 3.460764713265778, 31.146882419392, 0.00005672242602267426, 0.00005672242602267426),
    @("M2", "Col5a1", "Sdc3", "FAPs", 2, 0.6666666666666666, 0.1415573333333333, 0.424672, 0.0003000162501266517, 0.0003000162501266516, 3, 1, 18.11074233333333, 54.332227, 0.1400575643155068, 0.1400575643155068, 2.563708389393778, 23.073375504544, 0.0000420195452478107, 0.00004201954524781069),
    @("M2", "Col5a1", "Sdc3", "M2", 2, 0.6666666666666666, 0.1415573333333333, 0.424672, 0.0003000162501266517, 0.0003000162501266516, 3, 1, 80.684877, 242.054631, 0.6239682030546764, 0.6239682030546765, 11.421536028448, 102.793824256032, 0.0001872006004787292, 0.0001872006004787292),
    @("M2", "Col5a1", "Sdc3", "sCs", 2, 0.6666666666666666, 0.1415573333333333, 0.424672, 0.0003000162501266517, 0.0003000162501266516, 3, 1, 6.065862, 18.197586, 0.0469097202951384, 0.04690972029513841, 0.8586672490880001, 7.728005241792, 0.00001407367837743751, 0.00001407367837743751),
    @("sCs", "Col5a1", "Sdc3", "ECs", 3, 1, 54.08460233333333, 162.253807, 0.1146267678229633, 0.1146267678229633, 3, 1, 24.44779533333333, 73.343386, 0.1890645123346783, 0.1890645123346783, 1322.249288530056, 11900.2435967705, 0.02167185395894895, 0.02167185395894895),
    @("sCs", "Col5a1", "Sdc3", "FAPs", 3, 1, 54.08460233333333, 162.253807, 0.1146267678229633, 0.1146267678229633, 3, 1, 18.11074233333333, 54.332227, 0.1400575643155068, 0.1400575643155068, 979.5122970597988, 8815.610673538189, 0.01605434590664335, 0.01605434590664334),
    @("sCs", "Col5a1", "Sdc3", "M2", 3, 1, 54.08460233333333, 162.253807, 0.1146267678229633, 0.1146267678229633, 3, 1, 80.684877, 242.054631, 0.6239682030546764, 0.6239682030546765, 4363.809486858912, 39274.28538173022, 0.07152345834045999, 0.07152345834045999),
    @("sCs", "Col5a1", "Sdc3", "sCs", 3, 1, 54.08460233333333, 162.253807, 0.1146267678229633, 0.1146267678229633, 3, 1, 6.065862, 18.197586, 0.0469097202951384, 0.04690972029513841, 328.069734078878, 2952.627606709902, 0.005377109616910978, 0.005377109616910978)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowVals = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($r, $j + 1).Value = $rowVals[$j]
    }
}
